# "Ontbrekende geometrie in de views"
# Appends ", geometrie" to the VIEW ATTRIBUTEN (column F) values that were
# missing the geometrie column, and refreshes the related workbook/sheet
# view state the way Excel/LibreOffice do after such an edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three distinct "VIEW ATTRIBUTEN" texts (column F) that are missing
# the trailing "geometrie" attribute.
$old1 = "identificatie_lokaalid, bgt_type, plus_type"
$old2 = "identificatie_lokaalid, bgt_fysiekvoorkomen, plus_fysiekvoorkomen"
$old3 = "identificatie_lokaalid, bgt_functie, plus_functie"

$new1 = $old1 + ", geometrie"
$new2 = $old2 + ", geometrie"
$new3 = $old3 + ", geometrie"

# Data rows run from row 2 through row 132 (row 1 is the header).
for ($i = 2; $i -le 132; $i++) {
    $addr = "F" + $i
    $cell = $ws.Range($addr)
    $v = $cell.Value2

    if ($v -eq $old1) {
        $cell.Value = $new1
    } elseif ($v -eq $old2) {
        $cell.Value = $new2
    } elseif ($v -eq $old3) {
        $cell.Value = $new3
    }
}

# Editing the filtered range in LibreOffice/Excel produces a fresh,
# incrementally-numbered _FilterDatabase defined name for the sheet.
[void]$ws.Names.Add("_xlnm._FilterDatabase_0_0_0_0", "=Blad1!`$A`$1:`$J`$132")

# Rows 83 and 131 end up with a recalculated (non-default) row height.
$ws.Rows.Item(83).RowHeight = 13.8
$ws.Rows.Item(131).RowHeight = 13.8

# Leave the selection on the last touched cell, matching where the
# editor ended up after making the changes.
[void]$ws.Range("F132").Select()
